$wb = $excel.ActiveWorkbook

# "produkty" sheet (sheet1): replace "ryż" with "czekolada" in A8
$ws1 = $wb.Worksheets.Item("produkty")
$ws1.Range("A8").Value = "czekolada"

# Move the active selection to A9 to reflect the new cursor position
$ws1.Range("A9").Select()
